$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 39 (Z16_B01 / Kriminalität / Crime), which shifts rows 40-44 up by one.
$ws.Rows.Item(39).Delete()
